$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Color used for highlighted "new entry" cells (matches existing highlighted
# cells already present in the sheet, e.g. G24, H35, I49 -> fill rgb FF7B6993)
$highlightColor = 9660795

# Row 7 : Feb-2 (H) and Feb-4 (J) collections entered
$ws.Range("H7").Value = 620
$ws.Range("H7").Interior.Color = $highlightColor
$ws.Range("J7").Value = 7380
$ws.Range("J7").Interior.Color = $highlightColor

# Row 13 : Feb-4 (J) collection entered
$ws.Range("J13").Value = 5000

# Row 14 : Feb-4 (J) collection entered
$ws.Range("J14").Value = 2000

# Row 18 : Feb-4 (J) collection entered (highlighted)
$ws.Range("J18").Value = 5000
$ws.Range("J18").Interior.Color = $highlightColor

# Row 24 : Feb-4 (J) collection entered (highlighted)
$ws.Range("J24").Value = 5000
$ws.Range("J24").Interior.Color = $highlightColor

# Row 26 : Feb-4 (J) collection entered
$ws.Range("J26").Value = 1000

# Row 27 : Feb-4 (J) collection entered
$ws.Range("J27").Value = 2000

# Row 32 : Feb-4 (J) collection entered
$ws.Range("J32").Value = 1000

# Row 35 : Feb-4 (J) collection entered (highlighted)
$ws.Range("J35").Value = 3000
$ws.Range("J35").Interior.Color = $highlightColor

# Row 41 : Feb-4 (J) collection entered
$ws.Range("J41").Value = 3000

# Row 46 : Feb-4 (J) collection entered
$ws.Range("J46").Value = 2100

# Row 49 : Feb-4 (J) collection entered (highlighted)
$ws.Range("J49").Value = 3000
$ws.Range("J49").Interior.Color = $highlightColor

# Row 62 : Feb-4 (J) collection entered
$ws.Range("J62").Value = 4000

# Row 65 : Feb-4 (J) collection entered
$ws.Range("J65").Value = 1000

# Row 66 : Feb-4 (J) collection entered
$ws.Range("J66").Value = 5000

# Row 69 : Feb-4 (J) collection entered
$ws.Range("J69").Value = 1000

# Row 70 : Feb-4 (J) collection entered
$ws.Range("J70").Value = 1000

# Row 71 : Feb-4 (J) collection entered (highlighted)
$ws.Range("J71").Value = 6000
$ws.Range("J71").Interior.Color = $highlightColor

# Row 72 : Feb-4 (J) collection entered
$ws.Range("J72").Value = 1500

# Row 83 : Feb-3 (I) collection entered (highlighted)
$ws.Range("I83").Value = 5000
$ws.Range("I83").Interior.Color = $highlightColor

# Update the view: scroll/select so the last-entered cell is active,
# matching the saved cursor position from the edit session.
$ws.Range("J14").Select()
